# Weekly update: insert a new observation as the new most-recent row (row 23),
# pushing the existing rows 23-32 down to 24-33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 23 (shifts rows 23:32 -> 24:33)
$ws.Rows("23:23").Insert()

# Populate the new row 23 with the new weekly record
$ws.Range("A23").Value = 4
$ws.Range("B23").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value = "Los Lagos"
$ws.Range("D23").Value = 44518
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 300000000
$ws.Range("G23").Value = "Espárragos"
$ws.Range("H23").Value = "Sin especificar"
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 180
$ws.Range("K23").Value = 1600
$ws.Range("L23").Value = 1600
$ws.Range("M23").Value = 1600
$ws.Range("N23").Value = '$/kilo'
$ws.Range("O23").Value = "Provincia de Linares"
$ws.Range("P23").Value = 1600
$ws.Range("Q23").Value = 1
$ws.Range("R23").Value = "Hortaliza"
